$rows = @(
    @{ Row = 2; D = "29.993.61"; E = "  -0.46%  " },
    @{ Row = 3; D = "1.897.83"; E = "  -0.88%  " },
    @{ Row = 4; D = "1.000"; E = "  +0.07%  " },
    @{ Row = 5; D = "0.8296"; E = "  +4.85%  " },
    @{ Row = 6; D = "241.78"; E = "  -0.39%  " },
    @{ Row = 7; E = "  +0.03%  " },
    @{ Row = 8; D = "0.3271"; E = "  +2.77%  " },
    @{ Row = 9; D = "26.49"; E = "  +0.53%  " },
    @{ Row = 10; D = "0.07028"; E = "  +0.98%  " },
    @{ Row = 11; D = "0.08086"; E = "  +1.03%  " },
    @{ Row = 12; D = "0.7615"; E = "  +1.24%  " },
    @{ Row = 13; D = "1.914.33"; E = "  -0.03%  " },
    @{ Row = 14; D = "5.243"; E = "  +0.31%  " },
    @{ Row = 15; D = "92.13"; E = "  -1.47%  " },
    @{ Row = 16; D = "29.992.24"; E = "  -0.57%  " },
    @{ Row = 17; D = "14.08"; E = "  +0.10%  " },
    @{ Row = 18; E = "  -2.59%  " },
    @{ Row = 19; D = "243.53"; E = "  -2.18%  " },
    @{ Row = 20; D = "0.000007742"; E = "  -1.08%  " },
    @{ Row = 21; E = "  +0.11%  " },
    @{ Row = 22; D = "2.149.64"; E = "  -0.82%  " },
    @{ Row = 23; D = "1.000"; E = "  -0.05%  " },
    @{ Row = 24; D = "6.944"; E = "  -0.50%  " },
    @{ Row = 25; D = "0.1740"; E = "  +25.59%  " },
    @{ Row = 26; D = "9.250"; E = "  -0.77%  " },
    @{ Row = 27; D = "165.36"; E = "  -2.27%  " },
    @{ Row = 28; D = "18.89"; E = "  -0.53%  " },
    @{ Row = 29; D = "2.090"; E = "  +1.63%  " },
    @{ Row = 30; D = "1.361"; E = "  -2.16%  " },
    @{ Row = 31; D = "1.515"; E = "  -0.80%  " },
    @{ Row = 32; D = "0.05870"; E = "  +8.81%  " },
    @{ Row = 33; D = "4.273"; E = "  -2.13%  " },
    @{ Row = 34; D = "4.063"; E = "  -1.48%  " },
    @{ Row = 35; D = "1.264"; E = "  -0.21%  " },
    @{ Row = 36; D = "0.7302"; E = "  -1.25%  " },
    @{ Row = 37; D = "2.719"; E = "  -0.37%  " },
    @{ Row = 38; D = "0.01914"; E = "  -0.96%  " },
    @{ Row = 39; D = "2.775"; E = "  -0.56%  " },
    @{ Row = 40; D = "0.4435"; E = "  -0.57%  " },
    @{ Row = 41; D = "72.35"; E = "  -0.52%  " },
    @{ Row = 42; D = "5.848"; E = "  -5.46%  " },
    @{ Row = 43; D = "0.8518"; E = "  +2.07%  " },
    @{ Row = 44; E = "  +0.01%  " },
    @{ Row = 45; D = "1.896"; E = "  -0.53%  " },
    @{ Row = 46; D = "101.97"; E = "  +1.37%  " },
    @{ Row = 47; B = "EnergySwap"; C = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; D = "9.837"; E = "  -0.08%  " },
    @{ Row = 48; B = "Aptos"; C = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; D = "7.539"; E = "  -0.93%  " },
    @{ Row = 49; D = "989.57"; E = "  +2.48%  " },
    @{ Row = 50; D = "2.047.42"; E = "  -0.71%  " },
    @{ Row = 51; D = "1.517"; E = "  +0.53%  " }
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in $rows) {
    $rowNum = $r.Row

    if ($r.ContainsKey("B")) {
        $ws.Cells.Item($rowNum, 2).Value = $r.B
    }
    if ($r.ContainsKey("C")) {
        $ws.Cells.Item($rowNum, 3).Value = $r.C
    }
    if ($r.ContainsKey("D")) {
        $dCell = $ws.Cells.Item($rowNum, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $r.D
    }
    if ($r.ContainsKey("E")) {
        $ws.Cells.Item($rowNum, 5).Value = $r.E
    }
}
